# Add the "Form Tag" column (used by get_form_type in import_utils) to the
# capital-commitments sheet, defaulting every existing row to "Default".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in column P, directly after the existing "As Of" column (O).
$ws.Range("P1").Value = "Form Tag"

# Populate the new column for every data row already present in the sheet.
$lastRow = $ws.Cells(1048576, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 16).Value = "Default"
}

# Leave the new column selected, matching the post-edit cursor position.
$null = $ws.Range("P1").Select()
